$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# s2cDNAProtocol value changed from E7760 to E7420 for every data row (G2:G27)
$ws.Range("G2:G27").Value = "E7420"

# The protocol column picked up a slightly larger font along with the edit
$ws.Range("G2:G27").Font.Size = 11

# roboticS2Prep literal FALSE values were replaced with a live =FALSE() formula
$ws.Range("H2:H27").Formula = "=FALSE()"

# Selection left on the protocol column that was just edited
$ws.Range("G2:G27").Select()
